$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend header row with a 4th column "Email"
$ws.Range("D1").Value = "Email"
$ws.Range("D1").Font.Bold = $true

# New row 15 data (order matters for shared-string table insertion order)
$ws.Range("A15").Value = "Wordpress Blog"
$ws.Range("C15").Value = "Peoplespaceoc2"
$ws.Range("B15").Value = "womencoders-admin"
$ws.Range("D15").Value = "womencoders@gmail.com"

# Hyperlink on D15 pointing to the same mailto target used elsewhere for this address
$ws.Hyperlinks.Add($ws.Range("D15"), "mailto:womencoders@gmail.com")
$ws.Range("D15").Style = "Hyperlink"

# Update selection to match the new active cell
$ws.Range("B15").Select()
